$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value (serial 45189 = 2023-09-20) for every
# data row. The update bumps this date by one day (to 45190 = 2023-09-21) for all
# data rows, from row 2 through row 132.
$ws.Range("C2:C132").Value = 45190
